$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 'https://echovme.in/'
    3 = 'https://socialpanga.com/'
    4 = 'https://www.matebiz.com/'
    5 = 'https://clutch.co/in/web-designers'
    6 = 'https://www.matebiz.com/digital-marketing-service/'
    7 = 'https://webeesocial.com/'
    8 = 'https://seoexpertscompanyindia.com/'
    9 = 'https://www.designrush.com/agency/website-design-development/in'
    10 = 'https://www.designrush.com/agency/search-engine-optimization/in'
    11 = 'https://www.mandywebdesign.com/'
    12 = 'https://www.seotonic.com/'
    13 = 'https://www.ezrankings.com/seo-company-india.html'
    14 = 'https://www.designrush.com/agency/digital-agencies/in'
    15 = 'https://www.socialee.in/'
    16 = 'https://www.foduu.com/'
    17 = 'https://www.webhopers.com/top-web-designing-companies-in-india'
    18 = 'https://www.digitalsilk.com/'
    19 = 'https://www.webhopers.com/seo-company-india'
    20 = 'https://acodez.in/'
    21 = 'https://omrdigital.com/'
    22 = 'https://www.digidarts.com/'
    23 = 'https://www.whiteriversmedia.com/'
    24 = 'https://www.vocso.com/'
    25 = 'https://www.brandloom.com/'
    26 = 'https://www.seodiscovery.com/seo-company-india.php'
    27 = 'https://www.mumbaiwebdesign.in/'
    28 = 'https://colorwhistle.com/'
    29 = 'https://www.ezrankings.com/web-design-company-india.html'
    30 = 'https://florafountain.com/'
    31 = 'https://www.linkedin.com/pulse/top-10-digital-marketing-agencies-india-socialee-wldmf'
    32 = 'https://www.pixelcrayons.com/services/digital-marketing/seo'
    33 = 'https://thatware.co/seo-company-india/'
    34 = 'https://www.aoneseoservice.com/'
    35 = 'https://indiawebdesigns.in/'
    36 = 'https://www.aaravinfotech.com/web-design-services.php'
    37 = 'https://www.rankontechnologies.com/seo-services/'
    38 = 'https://www.techmagnate.com/'
    39 = 'https://seoserviceinindia.co.in/'
    40 = 'https://www.semrush.com/agencies/list/seo/india/'
    41 = 'https://www.ezrankings.com/'
    42 = 'https://www.rankingbyseo.com/'
    43 = 'https://pwskills.com/blog/digital-marketing-companies/'
    44 = 'https://noviindus.com/web-designing-company-india/'
    45 = 'https://www.orangemantra.com/services/search-engine-optimization/'
    46 = 'https://wefttechnologies.com/digital-marketing-company-in-india/'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
